$wb = $excel.ActiveWorkbook

$wsQ4 = $wb.Worksheets.Item("2021-Q4")
$wsOldTotal = $wb.Worksheets.Item("总计")     # will be renamed/repurposed -> "2022-Q1"

# --- Duplicate the current "总计" sheet (keeps sheetPr/pageMargins/etc. intact)
#     and place the copy at the end; this copy becomes the NEW "总计" summary
#     sheet once we update its data below. ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsOldTotal.Copy($null, $lastSheet)
$wsNewTotal = $wb.Worksheets.Item($wb.Worksheets.Count)

# Free up the "总计" name on the original sheet before claiming it on the copy
$wsOldTotal.Name = "2022-Q1"
$wsNewTotal.Name = "总计"

# --- Repurpose the original "总计" sheet as "2022-Q1", filled with the new
#     quarter's fund-holding data (same column layout as 2021-Q3 / 2021-Q4) ---
$wsOldTotal.Cells.Clear()

$wsQ4.Range("B1:H1").Copy($wsOldTotal.Range("B1:H1"))
$wsQ4.Range("A2").Copy($wsOldTotal.Range("A2"))

# Fund code / metrics are stored as TEXT (matching the 2021-Q3 / 2021-Q4
# sheets), not numbers, so pre-format the cells as Text before assigning
# numeric-looking strings -- otherwise the engine auto-coerces them to
# numbers. ClearFormats afterwards drops the now-unneeded "@" number format
# again (the stored value stays Text) so no stray style is left behind.
$wsOldTotal.Range("B2:G2").NumberFormat = "@"
$wsOldTotal.Range("B2").Value = "501081"
$wsOldTotal.Range("C2").Value = "中欧科创主题3年封闭运作灵活配置混合"
$wsOldTotal.Range("D2").Value = "16.69"
$wsOldTotal.Range("E2").Value = "77.88"
$wsOldTotal.Range("F2").Value = "4.13"
$wsOldTotal.Range("G2").Value = "0.6893"
$wsOldTotal.Range("B2:G2").ClearFormats()
$wsOldTotal.Range("H2").Value = 5

# --- Update the new "总计" sheet: insert a 2022-Q1 row at the top of the
#     data, pushing 2021-Q4 / 2021-Q3 down one row each ---
$wsNewTotal.Range("A3").Copy($wsNewTotal.Range("A4"))
$wsNewTotal.Range("A4").Value = 2
$wsNewTotal.Range("B4").Value = "2021-Q3"
$wsNewTotal.Range("C4").Value = 2
$wsNewTotal.Range("D4").Value = 0

$wsNewTotal.Range("A3").Value = 1
$wsNewTotal.Range("B3").Value = "2021-Q4"
$wsNewTotal.Range("C3").Value = 1
$wsNewTotal.Range("D3").Value = 0.84

$wsNewTotal.Range("A2").Value = 0
$wsNewTotal.Range("B2").Value = "2022-Q1"
$wsNewTotal.Range("C2").Value = 1
$wsNewTotal.Range("D2").Value = 0.6899999999999999
